$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "mixed"
$ws.Range("A12").Value = "mixed_form"
$ws.Range("A13").Value = "mixed_str"
$ws.Range("A14").Value = "discretiser"
$ws.Range("B14").Value = "total_credit_utilized,installment"

# Mirror the author's final cursor position after entering the new rows
# (scrolled down so row 7 is at the top, with the empty row below the
# new data selected).
$ws.Range("A15").Select()
